$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 121
$ws.Range("J2").Value = 22328
$ws.Range("J3").Value = 113030
$ws.Range("H4").Value = 7824
$ws.Range("J4").Value = 30017
$ws.Range("H5").Value = 5031
$ws.Range("I5").Value = 8331
$ws.Range("J5").Value = 53398

$ws.Range("F6").Select()
